# Apply pricing_update.xlsx template changes:
# - Add a new "Product Discounts" worksheet (Part Number / Discount)
# - Rework "Product Category Discounts" sheet to add Category Rank and
#   Pricing Category columns

$wb = $excel.ActiveWorkbook

# --- Sheet: "Product Category Discounts" (existing 3rd tab) ---
# Becomes: Product Category | Category Rank | Pricing Category | Discount
$catDiscounts = $wb.Worksheets.Item("Product Category Discounts")
$catDiscounts.Range("A1").Value = "Product Category"
$catDiscounts.Range("B1").Value = "Category Rank"
$catDiscounts.Range("C1").Value = "Pricing Category"
$catDiscounts.Range("D1").Value = "Discount"
$catDiscounts.Range("A1:D1").Style = "Normal"
$catDiscounts.Range("A1:D1").Font.Bold = $true
$catDiscounts.Range("B1:C1").ColumnWidth = 16.7109375

# --- New sheet: "Product Discounts" ---
$productDiscounts = $wb.Worksheets.Add()
$productDiscounts.Name = "Product Discounts"
$productDiscounts.Range("A1").Value = "Part Number"
$productDiscounts.Range("B1").Value = "Discount"
$productDiscounts.Range("A1:B1").Font.Bold = $true

# Move the new sheet to the end (after "Product Category Discounts")
$productDiscounts.Move($null, $wb.Worksheets.Item("Product Category Discounts"))
